$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-sort the two category pairs whose relative order flipped after
#     re-running the Russia (RU) computation: swap the row labels so the
#     numbers that already sit in rows 7/8 and 10/11 line up with their
#     (now reordered) category names.
$row7Label = $ws.Range("A7").Value2
$row8Label = $ws.Range("A8").Value2
$ws.Range("A7").Value = $row8Label
$ws.Range("A8").Value = $row7Label

$row10Label = $ws.Range("A10").Value2
$row11Label = $ws.Range("A11").Value2
$ws.Range("A10").Value = $row11Label
$ws.Range("A11").Value = $row10Label

# --- Re-run RU (Russia, column L) and refresh the "All" aggregate (column B);
#     this also carries fresh values for every other column on rows 7/8 and
#     10/11 because those rows were re-sorted above.

# Row 2
$ws.Range("B2").Value = 0.131625883609389
$ws.Range("L2").Value = 0.12261602368026

# Row 3
$ws.Range("B3").Value = 0.117391740788661
$ws.Range("L3").Value = 0.0780044195016593

# Row 4
$ws.Range("B4").Value = 0.115086501786906
$ws.Range("L4").Value = 0.233470317821513

# Row 5
$ws.Range("B5").Value = 0.0302250164353345
$ws.Range("L5").Value = 0.0400330890547246

# Row 6
$ws.Range("B6").Value = 0.0220629991092238
$ws.Range("L6").Value = 0.0263842405119762

# Row 7
$ws.Range("B7").Value = 0.020929216733633
$ws.Range("C7").Value = 0.0218347128817252
$ws.Range("D7").Value = 0.0218700451526134
$ws.Range("E7").Value = 0.0279920491282482
$ws.Range("F7").Value = 0.0133860546555798
$ws.Range("G7").Value = 0.0519893723564105
$ws.Range("H7").Value = 0.0281582658272028
$ws.Range("I7").Value = 0.0237144708035906
$ws.Range("J7").Value = 0.0584563140940413
$ws.Range("K7").Value = 0.0146902681742471
$ws.Range("L7").Value = 0.0109577271698444
$ws.Range("M7").Value = 0.0146942867926235
$ws.Range("N7").Value = 0.0232164381262679

# Row 8
$ws.Range("B8").Value = 0.0208245358406806
$ws.Range("C8").Value = 0.0108572546896661
$ws.Range("D8").Value = 0.0185212941677131
$ws.Range("E8").Value = 0.00699066316329541
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0.0219494164489717
$ws.Range("H8").Value = 0.0208606830072648
$ws.Range("I8").Value = 0.00976253072757774
$ws.Range("J8").Value = 0.0223141546648006
$ws.Range("K8").Value = 0.00488877169262998
$ws.Range("L8").Value = 0.0371790560030417
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0.0324960529700791

# Row 9
$ws.Range("B9").Value = 0.0171506810414515
$ws.Range("L9").Value = 0.0738200326131201

# Row 10
$ws.Range("B10").Value = 0.0100656430244752
$ws.Range("C10").Value = 0.00444800246628011
$ws.Range("D10").Value = 0.00335337435894486
$ws.Range("E10").Value = 0.000482625482625483
$ws.Range("F10").Value = 0.00855118235815349
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0.0026261196551152
$ws.Range("I10").Value = 0.00282289501256844
$ws.Range("J10").Value = 0.0018504561920746
$ws.Range("K10").Value = 0.00306997465436804
$ws.Range("L10").Value = 0.0148719786389175
$ws.Range("M10").Value = 0.00385418200790109
$ws.Range("N10").Value = 0.0129685698957936

# Row 11
$ws.Range("B11").Value = 0.00997359857741974
$ws.Range("C11").Value = 0.00728417020100042
$ws.Range("D11").Value = 0.0036971411514128
$ws.Range("E11").Value = 0.00430626252369435
$ws.Range("F11").Value = 0.00333954434204623
$ws.Range("G11").Value = 0.0119254166216743
$ws.Range("H11").Value = 0.00594200589182936
$ws.Range("I11").Value = 0.00646368141978556
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0.0121029296382471
$ws.Range("L11").Value = 0.0115254198961637
$ws.Range("M11").Value = 0.0219788440934882
$ws.Range("N11").Value = 0.00822570783875314

Write-Output "applied edits"